$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("F2").Value = 5
$ws.Range("I2").Value = 20

# Update the active selection on the sheet
$ws.Range("I14").Select()
